$d = $word.ActiveDocument

# Common run properties used throughout this document's body text.
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>'
$pPr = '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>'

function New-PkgXml([string]$bodyInner) {
    return @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyInner
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1) "COUNT().- " paragraph: split the lead run so "COUNT(" is wrapped in
#    proofErr gramStart/gramEnd marks, and drop the stray _GoBack bookmark
#    that used to sit at the end of this paragraph.
# ---------------------------------------------------------------------------
$countPara = $d.Paragraphs(4)
$countRng = $countPara.Range
if ($countRng.Text -notlike "COUNT*") { throw "Unexpected paragraph 4: $($countRng.Text)" }

$countBody = "<w:p>$pPr" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>COUNT(</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`">).- </w:t></w:r>" +
    "<w:r>$rPr<w:t>Sirve para ver la cantidad de registros que hay en una fila.</w:t></w:r>" +
    '</w:p>'
$countRng.InsertXML((New-PkgXml $countBody))

# ---------------------------------------------------------------------------
# 2) "SUM().- ..." paragraph: keep the existing SUM( / proofErr pair, split
#    the trailing text into three runs and fix "un fila o columna" -> "una
#    columna".
# ---------------------------------------------------------------------------
$sumPara = $d.Paragraphs(5)
$sumRng = $sumPara.Range
if ($sumRng.Text -notlike "SUM*") { throw "Unexpected paragraph 5: $($sumRng.Text)" }

$sumBody = "<w:p>$pPr" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>SUM(</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r>$rPr<w:t>).- Sirve para realizar una operación s</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">imple que es la suma en una </w:t></w:r>" +
    "<w:r>$rPr<w:t>columna.</w:t></w:r>" +
    '</w:p>'
$sumRng.InsertXML((New-PkgXml $sumBody))

# ---------------------------------------------------------------------------
# 3) "AVG.- ..." paragraph: wrap "AVG.-" in proofErr marks, split off the
#    following space into its own run, and fix "una fila o columna" -> "una
#    columna".
# ---------------------------------------------------------------------------
$avgPara = $d.Paragraphs(6)
$avgRng = $avgPara.Range
if ($avgRng.Text -notlike "AVG*") { throw "Unexpected paragraph 6: $($avgRng.Text)" }

$avgBody = "<w:p>$pPr" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>AVG.-</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> Sirve para el cálculo de </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> promedios ya sea de una columna.</w:t></w:r>" +
    '</w:p>'
$avgRng.InsertXML((New-PkgXml $avgBody))

# ---------------------------------------------------------------------------
# 4) "MIN.- ..." paragraph: wrap "MIN.-" in proofErr marks and split the
#    remainder into three runs, fixing "una fila o columna " -> "una columna ".
# ---------------------------------------------------------------------------
$minPara = $d.Paragraphs(7)
$minRng = $minPara.Range
if ($minRng.Text -notlike "MIN*") { throw "Unexpected paragraph 7: $($minRng.Text)" }

$minBody = "<w:p>$pPr" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>MIN.-</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> Utilizado para hall</w:t></w:r>" +
    "<w:r>$rPr<w:t>ar el valor mínimo de una</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> columna </w:t></w:r>" +
    '</w:p>'
$minRng.InsertXML((New-PkgXml $minBody))

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the COUNT paragraph (already removed
#    above, since that paragraph's XML was fully replaced) to the very end
#    of the document, right after the final period.
#
#    A Range collapsed exactly at Content.End seats incorrectly in this
#    runtime, so nudge past the boundary with a throwaway character, anchor
#    the bookmark there, then delete the throwaway character again.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$endPos = $d.Content.End
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$junk = $d.Range($endPos, $endPos + 1)
$junk.Delete()

Write-Host "Edit complete"
